# Horarios actualizados Linea 141 - 1323
# Applies the scheduled-data refresh (new scrape at 11:27:45) to all three
# sheets of the workbook: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 11:27:45"
$ws1.Cells.Item(3,1).Value = "Total filas: 165"

# --- Row pair swaps (A/C/D columns swap between the two rows; B/E unchanged) ---
$ws1.Cells.Item(56,1).Value = "07:35:06"
$ws1.Cells.Item(56,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(56,4).Value = 1

$ws1.Cells.Item(57,1).Value = "06:50:53"
$ws1.Cells.Item(57,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(57,4).Value = 46

$ws1.Cells.Item(67,3).Value = "10_OLMOS"
$ws1.Cells.Item(68,3).Value = "16_SANTA ANA"

$ws1.Cells.Item(73,1).Value = "06:50:53"
$ws1.Cells.Item(73,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(73,4).Value = 92

$ws1.Cells.Item(74,1).Value = "08:22:49"
$ws1.Cells.Item(74,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(74,4).Value = 0

$ws1.Cells.Item(85,1).Value = "08:54:22"
$ws1.Cells.Item(85,3).Value = "10_OLMOS"
$ws1.Cells.Item(85,4).Value = 0

$ws1.Cells.Item(86,1).Value = "08:22:49"
$ws1.Cells.Item(86,3).Value = "17_ROMERO"
$ws1.Cells.Item(86,4).Value = 32

$ws1.Cells.Item(100,1).Value = "07:35:06"
$ws1.Cells.Item(100,3).Value = "17_ROMERO"
$ws1.Cells.Item(100,4).Value = 108

$ws1.Cells.Item(102,1).Value = "08:54:22"
$ws1.Cells.Item(102,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(102,4).Value = 29

$ws1.Cells.Item(106,1).Value = "08:54:22"
$ws1.Cells.Item(106,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(106,4).Value = 40

$ws1.Cells.Item(107,1).Value = "08:22:49"
$ws1.Cells.Item(107,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(107,4).Value = 72

# --- New scrape batch (rows 142-158 refreshed, 159-170 newly appended) ---
$ws1.Cells.Item(142,1).Value = "11:27:45"
$ws1.Cells.Item(142,2).Value = "11:27"
$ws1.Cells.Item(142,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(142,4).Value = 0
$ws1.Cells.Item(142,5).Value = "LP1912"

$ws1.Cells.Item(143,1).Value = "11:27:45"
$ws1.Cells.Item(143,2).Value = "11:28"
$ws1.Cells.Item(143,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(143,4).Value = 1
$ws1.Cells.Item(143,5).Value = "LP1912"

$ws1.Cells.Item(144,1).Value = "11:27:45"
$ws1.Cells.Item(144,2).Value = "11:32"
$ws1.Cells.Item(144,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(144,4).Value = 5
$ws1.Cells.Item(144,5).Value = "LP1912"

$ws1.Cells.Item(145,1).Value = "11:27:45"
$ws1.Cells.Item(145,2).Value = "11:34"
$ws1.Cells.Item(145,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(145,4).Value = 7
$ws1.Cells.Item(145,5).Value = "LP1912"

$ws1.Cells.Item(146,1).Value = "11:27:45"
$ws1.Cells.Item(146,2).Value = "11:35"
$ws1.Cells.Item(146,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(146,4).Value = 8
$ws1.Cells.Item(146,5).Value = "LP1912"

$ws1.Cells.Item(147,1).Value = "11:27:45"
$ws1.Cells.Item(147,2).Value = "11:41"
$ws1.Cells.Item(147,3).Value = "17_ROMERO"
$ws1.Cells.Item(147,4).Value = 14
$ws1.Cells.Item(147,5).Value = "LP1912"

$ws1.Cells.Item(148,1).Value = "11:27:45"
$ws1.Cells.Item(148,2).Value = "11:43"
$ws1.Cells.Item(148,3).Value = "10_OLMOS"
$ws1.Cells.Item(148,4).Value = 16
$ws1.Cells.Item(148,5).Value = "LP1912"

$ws1.Cells.Item(149,1).Value = "11:27:45"
$ws1.Cells.Item(149,2).Value = "11:51"
$ws1.Cells.Item(149,3).Value = "215B_EL PATO"
$ws1.Cells.Item(149,4).Value = 24
$ws1.Cells.Item(149,5).Value = "LP1912"

$ws1.Cells.Item(150,1).Value = "11:27:45"
$ws1.Cells.Item(150,2).Value = "11:52"
$ws1.Cells.Item(150,3).Value = "15_ABASTO"
$ws1.Cells.Item(150,4).Value = 25
$ws1.Cells.Item(150,5).Value = "LP1912"

$ws1.Cells.Item(151,1).Value = "11:27:45"
$ws1.Cells.Item(151,2).Value = "11:59"
$ws1.Cells.Item(151,3).Value = "225_GOMEZ"
$ws1.Cells.Item(151,4).Value = 32
$ws1.Cells.Item(151,5).Value = "LP1912"

$ws1.Cells.Item(152,1).Value = "11:27:45"
$ws1.Cells.Item(152,2).Value = "12:02"
$ws1.Cells.Item(152,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(152,4).Value = 35
$ws1.Cells.Item(152,5).Value = "LP1912"

$ws1.Cells.Item(153,1).Value = "11:27:45"
$ws1.Cells.Item(153,2).Value = "12:04"
$ws1.Cells.Item(153,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(153,4).Value = 37
$ws1.Cells.Item(153,5).Value = "LP1912"

$ws1.Cells.Item(154,1).Value = "11:27:45"
$ws1.Cells.Item(154,2).Value = "12:06"
$ws1.Cells.Item(154,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(154,4).Value = 39
$ws1.Cells.Item(154,5).Value = "LP1912"

$ws1.Cells.Item(155,1).Value = "11:27:45"
$ws1.Cells.Item(155,2).Value = "12:06"
$ws1.Cells.Item(155,3).Value = "14_ABASTO"
$ws1.Cells.Item(155,4).Value = 39
$ws1.Cells.Item(155,5).Value = "LP1912"

# row 156 keeps its original "Hora_Scrap" (A156 = 10:45:47) per source data
$ws1.Cells.Item(156,1).Value = "10:45:47"
$ws1.Cells.Item(156,2).Value = "12:16"
$ws1.Cells.Item(156,3).Value = "17_ROMERO"
$ws1.Cells.Item(156,4).Value = 91
$ws1.Cells.Item(156,5).Value = "LP1912"

$ws1.Cells.Item(157,1).Value = "11:27:45"
$ws1.Cells.Item(157,2).Value = "12:20"
$ws1.Cells.Item(157,3).Value = "215A_EL PATO"
$ws1.Cells.Item(157,4).Value = 53
$ws1.Cells.Item(157,5).Value = "LP1912"

$ws1.Cells.Item(158,1).Value = "11:27:45"
$ws1.Cells.Item(158,2).Value = "12:20"
$ws1.Cells.Item(158,3).Value = "14_ABASTO"
$ws1.Cells.Item(158,4).Value = 53
$ws1.Cells.Item(158,5).Value = "LP1912"

$ws1.Cells.Item(159,1).Value = "11:27:45"
$ws1.Cells.Item(159,2).Value = "12:21"
$ws1.Cells.Item(159,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(159,4).Value = 54
$ws1.Cells.Item(159,5).Value = "LP1912"

$ws1.Cells.Item(160,1).Value = "11:27:45"
$ws1.Cells.Item(160,2).Value = "12:22"
$ws1.Cells.Item(160,3).Value = "10_OLMOS"
$ws1.Cells.Item(160,4).Value = 55
$ws1.Cells.Item(160,5).Value = "LP1912"

$ws1.Cells.Item(161,1).Value = "11:27:45"
$ws1.Cells.Item(161,2).Value = "12:36"
$ws1.Cells.Item(161,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(161,4).Value = 69
$ws1.Cells.Item(161,5).Value = "LP1912"

$ws1.Cells.Item(162,1).Value = "11:27:45"
$ws1.Cells.Item(162,2).Value = "12:38"
$ws1.Cells.Item(162,3).Value = "17_179 Y 38"
$ws1.Cells.Item(162,4).Value = 71
$ws1.Cells.Item(162,5).Value = "LP1912"

$ws1.Cells.Item(163,1).Value = "11:27:45"
$ws1.Cells.Item(163,2).Value = "12:41"
$ws1.Cells.Item(163,3).Value = "10_OLMOS"
$ws1.Cells.Item(163,4).Value = 74
$ws1.Cells.Item(163,5).Value = "LP1912"

$ws1.Cells.Item(164,1).Value = "11:27:45"
$ws1.Cells.Item(164,2).Value = "12:48"
$ws1.Cells.Item(164,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(164,4).Value = 81
$ws1.Cells.Item(164,5).Value = "LP1912"

$ws1.Cells.Item(165,1).Value = "11:27:45"
$ws1.Cells.Item(165,2).Value = "12:50"
$ws1.Cells.Item(165,3).Value = "15_ABASTO"
$ws1.Cells.Item(165,4).Value = 83
$ws1.Cells.Item(165,5).Value = "LP1912"

$ws1.Cells.Item(166,1).Value = "11:27:45"
$ws1.Cells.Item(166,2).Value = "12:58"
$ws1.Cells.Item(166,3).Value = "17_ROMERO"
$ws1.Cells.Item(166,4).Value = 91
$ws1.Cells.Item(166,5).Value = "LP1912"

$ws1.Cells.Item(167,1).Value = "11:27:45"
$ws1.Cells.Item(167,2).Value = "13:06"
$ws1.Cells.Item(167,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(167,4).Value = 99
$ws1.Cells.Item(167,5).Value = "LP1912"

$ws1.Cells.Item(168,1).Value = "11:27:45"
$ws1.Cells.Item(168,2).Value = "13:13"
$ws1.Cells.Item(168,3).Value = "215D_EL PATO"
$ws1.Cells.Item(168,4).Value = 106
$ws1.Cells.Item(168,5).Value = "LP1912"

$ws1.Cells.Item(169,1).Value = "11:27:45"
$ws1.Cells.Item(169,2).Value = "13:21"
$ws1.Cells.Item(169,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(169,4).Value = 114
$ws1.Cells.Item(169,5).Value = "LP1912"

$ws1.Cells.Item(170,1).Value = "11:27:45"
$ws1.Cells.Item(170,2).Value = "13:22"
$ws1.Cells.Item(170,3).Value = "10_OLMOS"
$ws1.Cells.Item(170,4).Value = 115
$ws1.Cells.Item(170,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 11:27:45"
$ws2.Cells.Item(3,1).Value = "Total filas: 23"

$ws2.Cells.Item(26,1).Value = "11:27:45"
$ws2.Cells.Item(26,4).Value = 24

$ws2.Cells.Item(27,1).Value = "11:27:45"
$ws2.Cells.Item(27,4).Value = 53

$ws2.Cells.Item(28,1).Value = "11:27:45"
$ws2.Cells.Item(28,2).Value = "13:13"
$ws2.Cells.Item(28,3).Value = "215D_EL PATO"
$ws2.Cells.Item(28,4).Value = 106
$ws2.Cells.Item(28,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 11:27:45"
$ws3.Cells.Item(3,1).Value = "Total filas: 23"

$ws3.Cells.Item(27,1).Value = "11:27:45"
$ws3.Cells.Item(27,4).Value = 37

$ws3.Cells.Item(28,1).Value = "11:27:45"
$ws3.Cells.Item(28,2).Value = "12:53"
$ws3.Cells.Item(28,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(28,4).Value = 86
$ws3.Cells.Item(28,5).Value = "L6203"
